# Update "想去人数" (interested-count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 484
$ws1.Range("F5").Value = 160
$ws1.Range("F7").Value = 658

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 484
$ws4.Range("F7").Value = 658
